$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Fri Mar 08 01:15:53 EST 2024"
$ws.Range("B4").Value = "Fri Mar 08 01:16:23 EST 2024"
